$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 220 (shifts existing rows 220:291 down to 221:292,
# matching Excel's "Insert Cells > Entire Row" behaviour / Rows.Insert).
$ws.Rows.Item(220).Insert()

# Populate the freshly inserted row 220 with the new weekly price record.
$ws.Cells.Item(220, 1).Value2  = 8
$ws.Cells.Item(220, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(220, 3).Value2  = "Coquimbo"
$ws.Cells.Item(220, 4).Value2  = 44876
$ws.Cells.Item(220, 5).Value2  = 4
$ws.Cells.Item(220, 6).Value2  = 100112031
$ws.Cells.Item(220, 7).Value2  = "Poroto verde"
$ws.Cells.Item(220, 8).Value2  = "Magnum"
$ws.Cells.Item(220, 9).Value2  = "Primera"
$ws.Cells.Item(220, 10).Value2 = 400
$ws.Cells.Item(220, 11).Value2 = 40000
$ws.Cells.Item(220, 12).Value2 = 41000
$ws.Cells.Item(220, 13).Value2 = 40500
$ws.Cells.Item(220, 14).Value2 = "`$/malla 25 kilos"
$ws.Cells.Item(220, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(220, 16).Value2 = 1620
$ws.Cells.Item(220, 17).Value2 = 25
$ws.Cells.Item(220, 18).Value2 = "Hortaliza"
